$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh prepends this week's two new price rows (Americana (o) /
# Inferno, "Primera" quality) for Comercializadora del Agro de Limarí - Ají,
# pushing the whole existing history down by two rows.
$ws.Rows("372:373").Insert()

$ws.Range("A372").Value = 2
$ws.Range("B372").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C372").Value = "Coquimbo"
$ws.Range("D372").Value = 45021
$ws.Range("E372").Value = 4
$ws.Range("F372").Value = 100112021
$ws.Range("G372").Value = "Ají"
$ws.Range("H372").Value = "Americana (o)"
$ws.Range("I372").Value = "Primera"
$ws.Range("J372").Value = 600
$ws.Range("K372").Value = 8000
$ws.Range("L372").Value = 9000
$ws.Range("M372").Value = 8500
$ws.Range("N372").Value = "$/caja 25 kilos"
$ws.Range("O372").Value = "Provincia de Limarí"
$ws.Range("P372").Value = 340
$ws.Range("Q372").Value = 25
$ws.Range("R372").Value = "Hortaliza"

$ws.Range("A373").Value = 2
$ws.Range("B373").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C373").Value = "Coquimbo"
$ws.Range("D373").Value = 45021
$ws.Range("E373").Value = 4
$ws.Range("F373").Value = 100112021
$ws.Range("G373").Value = "Ají"
$ws.Range("H373").Value = "Inferno"
$ws.Range("I373").Value = "Primera"
$ws.Range("J373").Value = 160
$ws.Range("K373").Value = 9000
$ws.Range("L373").Value = 10000
$ws.Range("M373").Value = 9500
$ws.Range("N373").Value = "$/caja 15 kilos"
$ws.Range("O373").Value = "Provincia de Limarí"
$ws.Range("P373").Value = 633
$ws.Range("Q373").Value = 15
$ws.Range("R373").Value = "Hortaliza"
